$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 29.223446
$ws.Range("H2").Value = 87.670338
$ws.Range("I2").Value = 0.0169041244192178
$ws.Range("J2").Value = 0.0169041244192178
$ws.Range("M2").Value = 8.820647333333334
$ws.Range("N2").Value = 26.461942
$ws.Range("O2").Value = 0.06415146660411865
$ws.Range("P2").Value = 0.06415146660411865
$ws.Range("Q2").Value = 257.7697110307107
$ws.Range("R2").Value = 2319.927399276396
$ws.Range("S2").Value = 0.001084424373151317
$ws.Range("T2").Value = 0.001084424373151317

# Row 3
$ws.Range("G3").Value = 29.223446
$ws.Range("H3").Value = 87.670338
$ws.Range("I3").Value = 0.0169041244192178
$ws.Range("J3").Value = 0.0169041244192178
$ws.Range("O3").Value = 0.3979101621202897
$ws.Range("P3").Value = 0.3979101621202898
$ws.Range("Q3").Value = 1598.85958864961
$ws.Range("R3").Value = 14389.73629784649
$ws.Range("S3").Value = 0.006726322888152505
$ws.Range("T3").Value = 0.006726322888152505

# Row 4
$ws.Range("G4").Value = 29.223446
$ws.Range("H4").Value = 87.670338
$ws.Range("I4").Value = 0.0169041244192178
$ws.Range("J4").Value = 0.0169041244192178
$ws.Range("M4").Value = 21.90816333333333
$ws.Range("N4").Value = 65.72449
$ws.Range("O4").Value = 0.1593353362087987
$ws.Range("P4").Value = 0.1593353362087987
$ws.Range("Q4").Value = 640.2320281308466
$ws.Range("R4").Value = 5762.088253177621
$ws.Range("S4").Value = 0.002693424347651432
$ws.Range("T4").Value = 0.002693424347651432

# Row 5
$ws.Range("G5").Value = 29.223446
$ws.Range("H5").Value = 87.670338
$ws.Range("I5").Value = 0.0169041244192178
$ws.Range("J5").Value = 0.0169041244192178
$ws.Range("M5").Value = 52.056859
$ws.Range("N5").Value = 156.170577
$ws.Range("O5").Value = 0.3786030350667928
$ws.Range("P5").Value = 0.3786030350667929
$ws.Range("Q5").Value = 1521.280807916114
$ws.Range("R5").Value = 13691.52727124502
$ws.Range("S5").Value = 0.006399952810262547
$ws.Range("T5").Value = 0.006399952810262548

# Row 6
$ws.Range("I6").Value = 0.9471112884046843
$ws.Range("J6").Value = 0.9471112884046842
$ws.Range("M6").Value = 8.820647333333334
$ws.Range("N6").Value = 26.461942
$ws.Range("O6").Value = 0.06415146660411865
$ws.Range("P6").Value = 0.06415146660411865
$ws.Range("Q6").Value = 14442.42819512425
$ws.Range("R6").Value = 129981.8537561183
$ws.Range("S6").Value = 0.06075857818847689
$ws.Range("T6").Value = 0.06075857818847688

# Row 7
$ws.Range("I7").Value = 0.9471112884046843
$ws.Range("J7").Value = 0.9471112884046842
$ws.Range("O7").Value = 0.3979101621202897
$ws.Range("P7").Value = 0.3979101621202898
$ws.Range("S7").Value = 0.3768652063150644
$ws.Range("T7").Value = 0.3768652063150644

# Row 8
$ws.Range("I8").Value = 0.9471112884046843
$ws.Range("J8").Value = 0.9471112884046842
$ws.Range("M8").Value = 21.90816333333333
$ws.Range("N8").Value = 65.72449
$ws.Range("O8").Value = 0.1593353362087987
$ws.Range("P8").Value = 0.1593353362087987
$ws.Range("Q8").Value = 35871.18539849274
$ws.Range("R8").Value = 322840.6685864347
$ws.Range("S8").Value = 0.1509082955651088
$ws.Range("T8").Value = 0.1509082955651088

# Row 9
$ws.Range("I9").Value = 0.9471112884046843
$ws.Range("J9").Value = 0.9471112884046842
$ws.Range("M9").Value = 52.056859
$ws.Range("N9").Value = 156.170577
$ws.Range("O9").Value = 0.3786030350667928
$ws.Range("P9").Value = 0.3786030350667929
$ws.Range("Q9").Value = 85234.95155849191
$ws.Range("R9").Value = 767114.5640264272
$ws.Range("S9").Value = 0.358579208336034
$ws.Range("T9").Value = 0.358579208336034

# Row 10
$ws.Range("G10").Value = 37.39212666666667
$ws.Range("H10").Value = 112.17638
$ws.Range("I10").Value = 0.02162924801792661
$ws.Range("J10").Value = 0.0216292480179266
$ws.Range("M10").Value = 8.820647333333334
$ws.Range("N10").Value = 26.461942
$ws.Range("O10").Value = 0.06415146660411865
$ws.Range("P10").Value = 0.06415146660411865
$ws.Range("Q10").Value = 329.8227623699956
$ws.Range("R10").Value = 2968.40486132996
$ws.Range("S10").Value = 0.001387547981894218
$ws.Range("T10").Value = 0.001387547981894218

# Row 11
$ws.Range("G11").Value = 37.39212666666667
$ws.Range("H11").Value = 112.17638
$ws.Range("I11").Value = 0.02162924801792661
$ws.Range("J11").Value = 0.0216292480179266
$ws.Range("O11").Value = 0.3979101621202897
$ws.Range("P11").Value = 0.3979101621202898
$ws.Range("Q11").Value = 2045.780646847767
$ws.Range("R11").Value = 18412.0258216299
$ws.Range("S11").Value = 0.008606497585353132
$ws.Range("T11").Value = 0.008606497585353132

# Row 12
$ws.Range("G12").Value = 37.39212666666667
$ws.Range("H12").Value = 112.17638
$ws.Range("I12").Value = 0.02162924801792661
$ws.Range("J12").Value = 0.0216292480179266
$ws.Range("M12").Value = 21.90816333333333
$ws.Range("N12").Value = 65.72449
$ws.Range("O12").Value = 0.1593353362087987
$ws.Range("P12").Value = 0.1593353362087987
$ws.Range("Q12").Value = 819.1928183940223
$ws.Range("R12").Value = 7372.735365546201
$ws.Range("S12").Value = 0.003446303504879828
$ws.Range("T12").Value = 0.003446303504879828

# Row 13
$ws.Range("G13").Value = 37.39212666666667
$ws.Range("H13").Value = 112.17638
$ws.Range("I13").Value = 0.02162924801792661
$ws.Range("J13").Value = 0.0216292480179266
$ws.Range("M13").Value = 52.056859
$ws.Range("N13").Value = 156.170577
$ws.Range("O13").Value = 0.3786030350667928
$ws.Range("P13").Value = 0.3786030350667929
$ws.Range("Q13").Value = 1946.516665596807
$ws.Range("R13").Value = 17518.64999037126
$ws.Range("S13").Value = 0.008188898945799427
$ws.Range("T13").Value = 0.008188898945799427

# Row 14
$ws.Range("G14").Value = 24.817167
$ws.Range("H14").Value = 74.45150100000001
$ws.Range("I14").Value = 0.01435533915817136
$ws.Range("J14").Value = 0.01435533915817136
$ws.Range("M14").Value = 8.820647333333334
$ws.Range("N14").Value = 26.461942
$ws.Range("O14").Value = 0.06415146660411865
$ws.Range("P14").Value = 0.06415146660411865
$ws.Range("Q14").Value = 218.903477919438
$ws.Range("R14").Value = 1970.131301274942
$ws.Range("S14").Value = 0.0009209160605962269
$ws.Range("T14").Value = 0.0009209160605962266

# Row 15
$ws.Range("G15").Value = 24.817167
$ws.Range("H15").Value = 74.45150100000001
$ws.Range("I15").Value = 0.01435533915817136
$ws.Range("J15").Value = 0.01435533915817136
$ws.Range("O15").Value = 0.3979101621202897
$ws.Range("P15").Value = 0.3979101621202898
$ws.Range("Q15").Value = 1357.785300921345
$ws.Range("R15").Value = 12220.06770829211
$ws.Range("S15").Value = 0.00571213533171971
$ws.Range("T15").Value = 0.00571213533171971

# Row 16
$ws.Range("G16").Value = 24.817167
$ws.Range("H16").Value = 74.45150100000001
$ws.Range("I16").Value = 0.01435533915817136
$ws.Range("J16").Value = 0.01435533915817136
$ws.Range("M16").Value = 21.90816333333333
$ws.Range("N16").Value = 65.72449
$ws.Range("O16").Value = 0.1593353362087987
$ws.Range("P16").Value = 0.1593353362087987
$ws.Range("Q16").Value = 543.69854810661
$ws.Range("R16").Value = 4893.286932959491
$ws.Range("S16").Value = 0.002287312791158567
$ws.Range("T16").Value = 0.002287312791158566

# Row 17
$ws.Range("G17").Value = 24.817167
$ws.Range("H17").Value = 74.45150100000001
$ws.Range("I17").Value = 0.01435533915817136
$ws.Range("J17").Value = 0.01435533915817136
$ws.Range("M17").Value = 52.056859
$ws.Range("N17").Value = 156.170577
$ws.Range("O17").Value = 0.3786030350667928
$ws.Range("P17").Value = 0.3786030350667929
$ws.Range("Q17").Value = 1291.903763298453
$ws.Range("R17").Value = 11627.13386968608
$ws.Range("S17").Value = 0.005434974974696856
$ws.Range("T17").Value = 0.005434974974696857

Write-Output "applied changes"